$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16502911309970775"
$ws1.Range("B2").Value = "go_stims-16502911309611669.csv"
$ws1.Range("B3").Value = "GNG_stims-16502911309805152.csv"
$ws1.Range("B4").Value = "go_stims-16502911309825144.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911309959836.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16502911328572485"
$ws2.Range("B2").Value = "TB-16502911323912656.csv"
$ws2.Range("B3").Value = "OB-16502911317623413.csv"
$ws2.Range("B4").Value = "TB-16502911328406.csv"
$ws2.Range("B5").Value = "OB-1650291131835157.csv"
$ws2.Range("B6").Value = "ZB-match_7-16502911312902994.csv"
$ws2.Range("B7").Value = "TB-16502911321693301.csv"
$ws2.Range("B8").Value = "ZB-match_5-16502911311286447.csv"
$ws2.Range("B9").Value = "OB-16502911316215563.csv"
$ws2.Range("B10").Value = "ZB-match_1-165029113153358.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16502911328572485"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16502911329198735"
$ws4.Range("B2").Value = "MM_stims-16502911328735456.csv"
$ws4.Range("B3").Value = "ZM_stims-1650291132860202.csv"
$ws4.Range("B4").Value = "MM_stims-16502911329040887.csv"
$ws4.Range("B5").Value = "ZM_stims-1650291132874549.csv"
$ws4.Range("B6").Value = "MM_stims-16502911329198735.csv"
$ws4.Range("B7").Value = "ZM_stims-165029113290509.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16502911329815648"
$ws5.Range("B2").Value = "vSAT_stims-16502911329517481.csv"
$ws5.Range("B3").Value = "SAT_stims-1650291132936884.csv"
$ws5.Range("B4").Value = "SAT_stims-16502911329236815.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650291132966182.csv"
